# Dagbækur fyrir alvöru karlmenn
# Swap the diary owner from Hildur Sif Thorarensen to Árni Víðir Jóhannesson
# and fill in the "Vika 6" week's time log (rows 28-30).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header: name / e-mail -------------------------------------------------
# (set e-mail first so the shared-string table order matches: Vika 8, e-mail, name)
$ws.Range("G2").Value = "avj2@hi.is"
$ws.Range("C2").Value = "Árni Víðir Jóhannesson"

# --- Vika 6 block (rows 26-32): fill in time entries ------------------------
# Hönnun (row 28) - þri
$ws.Range("D28").Value = 60
# Forritun (row 29) - mán / þri / mið
$ws.Range("C29").Value = 180
$ws.Range("D29").Value = 340
$ws.Range("E29").Value = 300
# Prófanir (row 30) - mið
$ws.Range("E30").Value = 60

# --- Cursor / selection ------------------------------------------------------
$ws.Range("M7").Select() | Out-Null
